# Weekly Work Report update:
#  - refresh the "Report Generated On" timestamp
#  - zero out the pricing column (audit/billing values pulled to $0)
#  - add two new line items (Point 30 GND-MD "Rem" and a new Point 31/Point 32
#    pair), which pushes the existing "Point 31 / CON-40-AAA-1-B" line down a
#    row and shifts the TOTAL row from row 27 to row 29
#  - update the summary counters (Total Billed Amount / Total Line Items)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Report generation timestamp
# ---------------------------------------------------------------------------
$ws.Range("D5").Value = "Report Generated On: 08/18/2025 09:49 PM"

# ---------------------------------------------------------------------------
# 2. Insert two fresh rows just before the old TOTAL row (row 27) so the
#    existing "Point 31" item (currently row 26) keeps its place, and the
#    TOTAL row slides from 27 down to 29. Excel auto-updates the A27:G27
#    merged TOTAL range to A29:G29 as part of the insert.
# ---------------------------------------------------------------------------
$ws.Rows("27:28").Insert()

# New row 27 must look like the "shaded" alternate-row style (same as rows
# 17/19/21/23/25) -- copy formatting from row 25. Copy only the used columns
# (A:H, matching the source rows exactly) rather than the whole row so the
# sheet's used-range/dimension isn't blown out to the full 16384 columns and
# no stray column-I cell gets created.
$ws.Range("A25:H25").Copy()
$ws.Range("A27:H27").PasteSpecial(-4122)   # xlPasteFormats

# New row 28 must look like the "plain" alternate-row style (same as rows
# 16/18/20/22/24/26) -- copy formatting from row 26.
$ws.Range("A26:H26").Copy()
$ws.Range("A28:H28").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3. Row 26 changes from "Point 31 / CON-40-AAA-1-B" to the new
#    "Point 30 / GND-MD" removal line (style/shading stays the same).
# ---------------------------------------------------------------------------
$ws.Range("A26").Value = "Point 30"
$ws.Range("B26").Value = "GND-MD"
$ws.Range("C26").Value = "Rem"
$ws.Range("D26").Value = "GND,Wire Mldg Only"
$ws.Range("E26").Value = "EA"
$ws.Range("F26").Value = 1
$ws.Range("H26").Value = 0

# ---------------------------------------------------------------------------
# 4. Row 27 (new) carries the item that used to live on row 26: the
#    "Point 31 / CON-40-AAA-1-B" install.
# ---------------------------------------------------------------------------
$ws.Range("A27").Value = "Point 31"
$ws.Range("B27").Value = "CON-40-AAA-1-B"
$ws.Range("C27").Value = "Inst"
$ws.Range("D27").Value = "CON,#4/0 AWG,Alum Alloy,One,Bare"
$ws.Range("E27").Value = "FT"
$ws.Range("F27").Value = 369
$ws.Range("H27").Value = 0

# ---------------------------------------------------------------------------
# 5. Row 28 (new) is a brand-new line item, Point 32.
# ---------------------------------------------------------------------------
$ws.Range("A28").Value = "Point 32"
$ws.Range("B28").Value = "ARM-8SF-GN-DL-C"
$ws.Range("C28").Value = "Inst"
$ws.Range("D28").Value = "ARM,8ftSgl.Fiberglass,Gain,DE Light,Corr"
$ws.Range("E28").Value = "EA"
$ws.Range("F28").Value = 0
$ws.Range("H28").Value = 0

# ---------------------------------------------------------------------------
# 6. Zero out every existing pricing cell (rows 16-25) -- all billed amounts
#    were reset to $0 in this audit pass.
# ---------------------------------------------------------------------------
foreach ($r in 16..25) {
    $ws.Cells.Item($r, 8).Value = 0
}

# ---------------------------------------------------------------------------
# 7. TOTAL row is now row 29 (moved automatically by the Insert above).
# ---------------------------------------------------------------------------
$ws.Range("H29").Value = 0

# ---------------------------------------------------------------------------
# 8. Summary box: Total Billed Amount -> $0, Total Line Items 11 -> 13.
# ---------------------------------------------------------------------------
$ws.Range("C8").Value = 0
$ws.Range("C9").Value = 13
